$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (site / email) appended below the existing header row.
$ws.Range("A2").Value = "https://spelevadores.com.br/"
$ws.Range("B2").Value = "contato@spelevadores.com.br;"
$ws.Range("A3").Value = "https://www.elevadoreskorman.com.br/empresas-elevadores-sao-paulo"
$ws.Range("B3").Value = "korman@elevadoreskorman.com.br;comercial@elevadoreskorman.com.br;vendas@elevadoreskorman.com.br;"

# Give the new rows the same look as the header (Arial 10, left aligned)
# but without the bold weight, by cloning the header's format and then
# switching bold off - this reuses the existing Arial/10 font definition
# instead of minting a pile of throwaway intermediate styles.
$ws.Range("A1:B1").Copy()
$dataRng = $ws.Range("A2:B3")
$dataRng.PasteSpecial(-4122)
$dataRng.Font.Bold = $false
